$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")
$r = $ws.Range("H5")
Write-Output ($r | Get-Member)
